$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update financial data values for rows 2-6 (columns D through AJ)
$ws.Range("D2").Value = 1150
$ws.Range("E2").Value = 75
$ws.Range("F2").Value = 75
$ws.Range("G2").Value = 68
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 57
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 1334
$ws.Range("L2").Value = 438
$ws.Range("M2").Value = 896
$ws.Range("N2").Value = 871
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 252
$ws.Range("Q2").Value = 49
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = -39
$ws.Range("T2").Value = 12
$ws.Range("U2").Value = 37
$ws.Range("V2").Value = 284
$ws.Range("W2").Value = 6.52
$ws.Range("X2").Value = 5.19
$ws.Range("Y2").Value = 6.76
$ws.Range("Z2").Value = 4.43
$ws.Range("AA2").Value = 48.92
$ws.Range("AB2").Value = 242.72
$ws.Range("AC2").Value = 113
$ws.Range("AD2").Value = 13.45
$ws.Range("AE2").Value = 1726
$ws.Range("AF2").Value = 0.88
$ws.Range("AG2").Value = 15
$ws.Range("AH2").Value = 0.99
$ws.Range("AI2").Value = 13.31
$ws.Range("AJ2").Value = 50450590
$ws.Range("D3").Value = 1018
$ws.Range("E3").Value = 45
$ws.Range("F3").Value = 45
$ws.Range("G3").Value = 29
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 1192
$ws.Range("L3").Value = 301
$ws.Range("M3").Value = 892
$ws.Range("N3").Value = 876
$ws.Range("O3").Value = 15
$ws.Range("P3").Value = 252
$ws.Range("Q3").Value = 164
$ws.Range("R3").Value = -22
$ws.Range("S3").Value = -137
$ws.Range("T3").Value = 28
$ws.Range("U3").Value = 135
$ws.Range("V3").Value = 161
$ws.Range("W3").Value = 4.45
$ws.Range("X3").Value = 1.26
$ws.Range("Y3").Value = 1.23
$ws.Range("Z3").Value = 1.02
$ws.Range("AA3").Value = 33.74
$ws.Range("AB3").Value = 243.32
$ws.Range("AC3").Value = 21
$ws.Range("AD3").Value = 128.13
$ws.Range("AE3").Value = 1737
$ws.Range("AF3").Value = 1.58
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 50450590
$ws.Range("D4").Value = 888
$ws.Range("E4").Value = 39
$ws.Range("F4").Value = 39
$ws.Range("G4").Value = 41
$ws.Range("H4").Value = 29
$ws.Range("I4").Value = 27
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 1086
$ws.Range("L4").Value = 173
$ws.Range("M4").Value = 913
$ws.Range("N4").Value = 913
$ws.Range("P4").Value = 270
$ws.Range("Q4").Value = 128
$ws.Range("R4").Value = -38
$ws.Range("S4").Value = -120
$ws.Range("T4").Value = 18
$ws.Range("U4").Value = 110
$ws.Range("V4").Value = 41
$ws.Range("W4").Value = 4.4
$ws.Range("X4").Value = 3.27
$ws.Range("Y4").Value = 2.98
$ws.Range("Z4").Value = 2.55
$ws.Range("AA4").Value = 18.99
$ws.Range("AB4").Value = 236.75
$ws.Range("AC4").Value = 51
$ws.Range("AD4").Value = 39.26
$ws.Range("AE4").Value = 1774
$ws.Range("AF4").Value = 1.14
$ws.Range("AG4").Value = 10
$ws.Range("AH4").Value = 0.5
$ws.Range("AI4").Value = 19.31
$ws.Range("AJ4").Value = 53985163
$ws.Range("D5").Value = 1019
$ws.Range("E5").Value = 33
$ws.Range("F5").Value = 33
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = 29
$ws.Range("I5").Value = 29
$ws.Range("K5").Value = 1105
$ws.Range("L5").Value = 181
$ws.Range("M5").Value = 924
$ws.Range("N5").Value = 924
$ws.Range("P5").Value = 270
$ws.Range("Q5").Value = 66
$ws.Range("R5").Value = 9
$ws.Range("S5").Value = -21
$ws.Range("T5").Value = 16
$ws.Range("U5").Value = 50
$ws.Range("V5").Value = 24
$ws.Range("W5").Value = 3.24
$ws.Range("X5").Value = 2.84
$ws.Range("Y5").Value = 3.15
$ws.Range("Z5").Value = 2.64
$ws.Range("AA5").Value = 19.54
$ws.Range("AB5").Value = 244.7
$ws.Range("AC5").Value = 54
$ws.Range("AD5").Value = 26.3
$ws.Range("AE5").Value = 1796
$ws.Range("AF5").Value = 0.78
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 53985163
$ws.Range("D6").Value = 846
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 88
$ws.Range("H6").Value = 67
$ws.Range("I6").Value = 67
$ws.Range("K6").Value = 2049
$ws.Range("L6").Value = 1060
$ws.Range("M6").Value = 990
$ws.Range("N6").Value = 990
$ws.Range("P6").Value = 270
$ws.Range("Q6").Value = 23
$ws.Range("R6").Value = -661
$ws.Range("S6").Value = 645
$ws.Range("T6").Value = 9
$ws.Range("U6").Value = 15
$ws.Range("V6").Value = 731
$ws.Range("W6").Value = 0.5
$ws.Range("X6").Value = 7.86
$ws.Range("Y6").Value = 6.95
$ws.Range("Z6").Value = 4.22
$ws.Range("AA6").Value = 107.03
$ws.Range("AB6").Value = 268.25
$ws.Range("AC6").Value = 123
$ws.Range("AD6").Value = 13.39
$ws.Range("AE6").Value = 1924
$ws.Range("AF6").Value = 0.86
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 53985163

# Remove cells that no longer have data (rows 4-6)
$ws.Range("O4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# Rows 7-9: remove all financial data, keep only A/B/C identifying columns
$ws.Range("D7:AJ9").ClearContents()
